# Dayton Ohio - Rethinking the Employer Portfolio for Population Growth - FINAL.pptx
#
# Slide 6 ("The Conclusion"), Content Placeholder 2:
#   - Paragraph 2: reword the sentence about service providers and split it
#     into three runs ("...competition " / "and stifling " / "innovation...").
#   - Paragraph 3: merge its three runs back into a single run (text itself
#     stays the same).
#
# Both edits are done by keeping one trailing "anchor" character from the
# run being replaced (so the new text can inherit its run formatting via
# InsertBefore), inserting the desired text in front of it, and finally
# deleting that now-redundant anchor character.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# ---------------------------------------------------------------------------
# Paragraph 2 -> three runs
# ---------------------------------------------------------------------------
$para2 = $tr.Paragraphs(2, 1)

# A paragraph's Text/Length include the trailing paragraph-mark characters
# (vertical-tab + carriage-return) for every paragraph except the very last
# one in the text frame, so strip those two off to get the real text length.
$breakMarkerLen = $para2.Length - ($tr.Characters($para2.Start, $para2.Length).Text.TrimEnd("`v", "`r").Length)
$textLen2 = $para2.Length - $breakMarkerLen

# Delete all but the final character of the old sentence; that last
# character stays behind as a formatting anchor.
$oldPart2 = $tr.Characters($para2.Start, $textLen2 - 1)
$oldPart2.Delete()

$para2b = $tr.Paragraphs(2, 1)

# Insert the three replacement runs (in reverse order) right before the
# anchor character; each one inherits the anchor run's formatting.
$null = $para2b.InsertBefore("innovation in other sectors.")
$null = $para2b.InsertBefore("and stifling ")
$null = $para2b.InsertBefore("Dayton’s economic portfolio has an overabundance of service providers, which could be hurting competition ")

# Remove the leftover anchor character (the old sentence's final period).
$para2c = $tr.Paragraphs(2, 1)
$anchorPos2 = $para2c.Start + $para2c.Length - $breakMarkerLen - 1
$tr.Characters($anchorPos2, 1).Delete()

# ---------------------------------------------------------------------------
# Paragraph 3 -> merge back into a single run
# ---------------------------------------------------------------------------
$para3 = $tr.Paragraphs(3, 1)

$breakMarkerLen3 = $para3.Length - ($tr.Characters($para3.Start, $para3.Length).Text.TrimEnd("`v", "`r").Length)
$textLen3 = $para3.Length - $breakMarkerLen3

$oldPart3 = $tr.Characters($para3.Start, $textLen3 - 1)
$oldPart3.Delete()

$para3b = $tr.Paragraphs(3, 1)
$null = $para3b.InsertBefore("Dayton policymakers should consider implementing negative tax incentives to reduce the proportion of service providers, then allocate the tax revenues into promoting and building residential, cultural, recreational, commercial, and transportation infrastructure.")

$para3c = $tr.Paragraphs(3, 1)
$anchorPos3 = $para3c.Start + $para3c.Length - $breakMarkerLen3 - 1
$tr.Characters($anchorPos3, 1).Delete()
